$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Fecha 16-nov-2022 -> 28-dic-2022 (Segunda/100 -> Primera/80, price cols updated)
$ws.Range("D2").Value = 44923
$ws.Range("L2").Value = "Primera"
$ws.Range("M2").Value = 80
$ws.Range("N2").Value = 7500
$ws.Range("O2").Value = 8000
$ws.Range("P2").Value = 7625
$ws.Range("S2").Value = 7625

# Row 3: Fecha 09-nov-2022 -> 28-dic-2022 (Volumen 200 -> 80, prices updated)
$ws.Range("D3").Value = 44923
$ws.Range("M3").Value = 80
$ws.Range("P3").Value = 7625
$ws.Range("S3").Value = 7625

# Row 4: Fecha 28-dic-2022 -> 16-nov-2022 (Primera/80 -> Segunda/100, prices updated)
$ws.Range("D4").Value = 44881
$ws.Range("L4").Value = "Segunda"
$ws.Range("M4").Value = 100
$ws.Range("N4").Value = 11250
$ws.Range("O4").Value = 11250
$ws.Range("P4").Value = 11250
$ws.Range("S4").Value = 11250

# Row 5: Fecha 28-dic-2022 -> 09-nov-2022 (Volumen 80 -> 200, prices updated)
$ws.Range("D5").Value = 44874
$ws.Range("M5").Value = 200
$ws.Range("P5").Value = 7750
$ws.Range("S5").Value = 7750
